$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly-reported day's data on row 47 (date serial 46006 = 2025-12-15)
$ws.Range("A47").Value = 46006
$ws.Range("B47").Value = 772
$ws.Range("C47").Value = 23
$ws.Range("D47").Value = 749

# Move the active selection down to the row just filled in, matching the
# author's last on-screen selection after entering the new data.
$ws.Range("A47:D47").Select()
